$p = $ppt.ActivePresentation

# ---------------------------------------------------------------------------
# 1) Handout master / notes master "datetimeFigureOut" fields: 9/1/2021 -> 3/25/2022
#    (best effort -- some runtimes do not allow editing master placeholder
#    field text; guarded so a failure here does not abort the rest of the
#    script)
# ---------------------------------------------------------------------------
try {
    $hm = $p.HandoutMaster
    $hmDate = $hm.Shapes.Item(2)
    $hmDate.TextFrame.TextRange.Text = "3/25/2022"
} catch {
    Write-Output "HandoutMaster date field could not be updated: $_"
}

try {
    $nm = $p.NotesMaster
    $nmDate = $nm.Shapes.Item(2)
    $nmDate.TextFrame.TextRange.Text = "3/25/2022"
} catch {
    Write-Output "NotesMaster date field could not be updated: $_"
}

# ---------------------------------------------------------------------------
# 2) Slide 1 "License and Citation" paragraph: update citation sentence and
#    split the DOI hyperlink run so it reads 10.6084/m9.figshare.19416767
# ---------------------------------------------------------------------------
$s = $p.Slides.Item(1)
$sh = $s.Shapes.Item(2)
$tr = $sh.TextFrame.TextRange

$oldSentence = "The requested citation the overall tutorial is: David E. Bernholdt, Anshu Dubey, Patricia A. Grubel, Rinku K. Gupta, and Gregory R. Watson, Better Scientific Software tutorial, in the International Conference for High-Performance Computing, Networking, Storage, and Analysis (SC21), St. Louis, MO, USA and online, 2021. DOI: "
$newSentence = "The requested citation the overall tutorial is: David E. Bernholdt, Patricia A. Grubel, Rinku K. Gupta, and David M. Rogers, Better Scientific Software tutorial, in Improving Scientific Software conference, online, 2022. DOI: "

$sentenceStart = $tr.Text.IndexOf($oldSentence) + 1
$sentenceRange = $tr.Characters($sentenceStart, $oldSentence.Length)
$sentenceRange.Text = $newSentence

# Re-fetch the text range/offsets after the sentence-length change, then trim
# the trailing ".16556628" off the DOI hyperlink run -- PowerPoint splits the
# touched run off from the untouched "10.6084/m9.figshare" prefix, so the
# hyperlink/bold formatting on both halves is preserved automatically.
$tr2 = $sh.TextFrame.TextRange
$oldSuffix = ".16556628"
$newSuffix = ".19416767"
$suffixStart = $tr2.Text.IndexOf($oldSuffix) + 1
$suffixRange = $tr2.Characters($suffixStart, $oldSuffix.Length)
$suffixRange.Text = $newSuffix
